$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 9,20
$data[0,0] = "FAPs"
$data[0,1] = "ECs"
$data[0,2] = "Eng"
$data[0,3] = "FAPs"
$data[0,4] = 2
$data[0,5] = 0.6666666666666666
$data[0,6] = 9.163165666666666
$data[0,7] = 27.489497
$data[0,8] = 0.2800251397703982
$data[0,9] = 0.2800251397703982
$data[0,10] = 3
$data[0,11] = 1
$data[0,12] = 114.155417
$data[0,13] = 342.466251
$data[0,14] = 0.6835107367845005
$data[0,15] = 0.6835107367845005
$data[0,16] = 1046.024997718416
$data[0,17] = 9414.224979465747
$data[0,18] = 0.1914001896026476
$data[0,19] = 0.1914001896026476
$data[1,0] = "FAPs"
$data[1,1] = "ECs"
$data[1,2] = "Eng"
$data[1,3] = "sCs"
$data[1,4] = 2
$data[1,5] = 0.6666666666666666
$data[1,6] = 9.163165666666666
$data[1,7] = 27.489497
$data[1,8] = 0.2800251397703982
$data[1,9] = 0.2800251397703982
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 35.924535
$data[1,13] = 107.773605
$data[1,14] = 0.2150997826628812
$data[1,15] = 0.2150997826628812
$data[1,16] = 329.182465702965
$data[1,17] = 2962.642191326685
$data[1,18] = 0.06023334670475557
$data[1,19] = 0.06023334670475559
$data[2,0] = "FAPs"
$data[2,1] = "ECs"
$data[2,2] = "Eng"
$data[2,3] = "Bmp2"
$data[2,4] = 2
$data[2,5] = 0.6666666666666666
$data[2,6] = 9.163165666666666
$data[2,7] = 27.489497
$data[2,8] = 0.2800251397703982
$data[2,9] = 0.2800251397703982
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 16.93339666666667
$data[2,13] = 50.80019
$data[2,14] = 0.1013894805526183
$data[2,15] = 0.1013894805526183
$data[2,16] = 155.1635189560478
$data[2,17] = 1396.47167060443
$data[2,18] = 0.02839160346299502
$data[2,19] = 0.02839160346299503
$data[3,0] = "sCs"
$data[3,1] = "ECs"
$data[3,2] = "Eng"
$data[3,3] = "FAPs"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 12.06704733333333
$data[3,7] = 36.201142
$data[3,8] = 0.3687673822623249
$data[3,9] = 0.3687673822623249
$data[3,10] = 3
$data[3,11] = 1
$data[3,12] = 114.155417
$data[3,13] = 342.466251
$data[3,14] = 0.6835107367845005
$data[3,15] = 0.6835107367845005
$data[3,16] = 1377.518820295405
$data[3,17] = 12397.66938265864
$data[3,18] = 0.2520564651522132
$data[3,19] = 0.2520564651522132
$data[4,0] = "sCs"
$data[4,1] = "ECs"
$data[4,2] = "Eng"
$data[4,3] = "sCs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 12.06704733333333
$data[4,7] = 36.201142
$data[4,8] = 0.3687673822623249
$data[4,9] = 0.3687673822623249
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 35.924535
$data[4,13] = 107.773605
$data[4,14] = 0.2150997826628812
$data[4,15] = 0.2150997826628812
$data[4,16] = 433.5030642729899
$data[4,17] = 3901.52757845691
$data[4,18] = 0.0793217837777857
$data[4,19] = 0.07932178377778572
$data[5,0] = "sCs"
$data[5,1] = "ECs"
$data[5,2] = "Eng"
$data[5,3] = "Bmp2"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 12.06704733333333
$data[5,7] = 36.201142
$data[5,8] = 0.3687673822623249
$data[5,9] = 0.3687673822623249
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 16.93339666666667
$data[5,13] = 50.80019
$data[5,14] = 0.1013894805526183
$data[5,15] = 0.1013894805526183
$data[5,16] = 204.3360990907755
$data[5,17] = 1839.02489181698
$data[5,18] = 0.03738913333232596
$data[5,19] = 0.03738913333232596
$data[6,0] = "Bmp2"
$data[6,1] = "ECs"
$data[6,2] = "Eng"
$data[6,3] = "FAPs"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 11.49244066666667
$data[6,7] = 34.477322
$data[6,8] = 0.3512074779672769
$data[6,9] = 0.351207477967277
$data[6,10] = 3
$data[6,11] = 1
$data[6,12] = 114.155417
$data[6,13] = 342.466251
$data[6,14] = 0.6835107367845005
$data[6,15] = 0.6835107367845005
$data[6,16] = 1311.924356651091
$data[6,17] = 11807.31920985982
$data[6,18] = 0.2400540820296397
$data[6,19] = 0.2400540820296397
$data[7,0] = "Bmp2"
$data[7,1] = "ECs"
$data[7,2] = "Eng"
$data[7,3] = "sCs"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 11.49244066666667
$data[7,7] = 34.477322
$data[7,8] = 0.3512074779672769
$data[7,9] = 0.351207477967277
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 35.924535
$data[7,13] = 107.773605
$data[7,14] = 0.2150997826628812
$data[7,15] = 0.2150997826628812
$data[7,16] = 412.86058696509
$data[7,17] = 3715.74528268581
$data[7,18] = 0.0755446521803399
$data[7,19] = 0.07554465218033991
$data[8,0] = "Bmp2"
$data[8,1] = "ECs"
$data[8,2] = "Eng"
$data[8,3] = "Bmp2"
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 11.49244066666667
$data[8,7] = 34.477322
$data[8,8] = 0.3512074779672769
$data[8,9] = 0.351207477967277
$data[8,10] = 3
$data[8,11] = 1
$data[8,12] = 16.93339666666667
$data[8,13] = 50.80019
$data[8,14] = 0.1013894805526183
$data[8,15] = 0.1013894805526183
$data[8,16] = 194.6060564767978
$data[8,17] = 1751.45450829118
$data[8,18] = 0.03560874375729736
$data[8,19] = 0.03560874375729736

$ws.Range("A2:T10").Value = $data

